$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$timeTaken = @(
    "2021-10-05 14:20:18.400064",
    "2021-10-05 14:20:18.400073",
    "2021-10-05 14:20:18.400078",
    "2021-10-05 14:20:18.400083",
    "2021-10-05 14:20:18.400088",
    "2021-10-05 14:20:18.400092",
    "2021-10-05 14:20:18.400096",
    "2021-10-05 14:20:18.400100",
    "2021-10-05 14:20:18.400105",
    "2021-10-05 14:20:18.400109",
    "2021-10-05 14:20:18.400111",
    "2021-10-05 14:20:18.400114",
    "2021-10-05 14:20:18.400116",
    "2021-10-05 14:20:18.400118",
    "2021-10-05 14:20:18.400121",
    "2021-10-05 14:20:18.400123",
    "2021-10-05 14:20:18.400126",
    "2021-10-05 14:20:18.400129",
    "2021-10-05 14:20:18.400131",
    "2021-10-05 14:20:18.400134",
    "2021-10-05 14:20:18.400136",
    "2021-10-05 14:20:18.400138",
    "2021-10-05 14:20:18.400141",
    "2021-10-05 14:20:18.400143",
    "2021-10-05 14:20:18.400146",
    "2021-10-05 14:20:18.400149",
    "2021-10-05 14:20:18.400151",
    "2021-10-05 14:20:18.400153",
    "2021-10-05 14:20:18.400156",
    "2021-10-05 14:20:18.400158",
    "2021-10-05 14:20:18.400161",
    "2021-10-05 14:20:18.400163",
    "2021-10-05 14:20:18.400166",
    "2021-10-05 14:20:18.400169",
    "2021-10-05 14:20:18.400171",
    "2021-10-05 14:20:18.400173",
    "2021-10-05 14:20:18.400176",
    "2021-10-05 14:20:18.400178",
    "2021-10-05 14:20:18.400181",
    "2021-10-05 14:20:18.400183",
    "2021-10-05 14:20:18.400186",
    "2021-10-05 14:20:18.400189",
    "2021-10-05 14:20:18.400191",
    "2021-10-05 14:20:18.400194",
    "2021-10-05 14:20:18.400196",
    "2021-10-05 14:20:18.400199",
    "2021-10-05 14:20:18.400201",
    "2021-10-05 14:20:18.400204",
    "2021-10-05 14:20:18.400206",
    "2021-10-05 14:20:18.400208",
    "2021-10-05 14:20:18.400211",
    "2021-10-05 14:20:18.400213",
    "2021-10-05 14:20:18.400216",
    "2021-10-05 14:20:18.400219",
    "2021-10-05 14:20:18.400221",
    "2021-10-05 14:20:18.400224",
    "2021-10-05 14:20:18.400226",
    "2021-10-05 14:20:18.400229",
    "2021-10-05 14:20:18.400231",
    "2021-10-05 14:20:18.400234",
    "2021-10-05 14:20:18.400236",
    "2021-10-05 14:20:18.400238",
    "2021-10-05 14:20:18.400241",
    "2021-10-05 14:20:18.400244",
    "2021-10-05 14:20:18.400247",
    "2021-10-05 14:20:18.400250",
    "2021-10-05 14:20:18.400253",
    "2021-10-05 14:20:18.400255",
    "2021-10-05 14:20:18.400257",
    "2021-10-05 14:20:18.400260",
    "2021-10-05 14:20:18.400262",
    "2021-10-05 14:20:18.400265",
    "2021-10-05 14:20:18.400267",
    "2021-10-05 14:20:18.400270",
    "2021-10-05 14:20:18.400272",
    "2021-10-05 14:20:18.400275",
    "2021-10-05 14:20:18.400280",
    "2021-10-05 14:20:18.400283",
    "2021-10-05 14:20:18.400285",
    "2021-10-05 14:20:18.400288",
    "2021-10-05 14:20:18.400290",
    "2021-10-05 14:20:18.400293",
    "2021-10-05 14:20:18.400296",
    "2021-10-05 14:20:18.400298",
    "2021-10-05 14:20:18.400300",
    "2021-10-05 14:20:18.400303",
    "2021-10-05 14:20:18.400305",
    "2021-10-05 14:20:18.400308",
    "2021-10-05 14:20:18.400310",
    "2021-10-05 14:20:18.400313",
    "2021-10-05 14:20:18.400315",
    "2021-10-05 14:20:18.400317",
    "2021-10-05 14:20:18.400321",
    "2021-10-05 14:20:18.400324",
    "2021-10-05 14:20:18.400326",
    "2021-10-05 14:20:18.400329",
    "2021-10-05 14:20:18.400331",
    "2021-10-05 14:20:18.400334",
    "2021-10-05 14:20:18.400336",
    "2021-10-05 14:20:18.400339",
    "2021-10-05 14:20:18.400341",
    "2021-10-05 14:20:18.400344",
    "2021-10-05 14:20:18.400346",
    "2021-10-05 14:20:18.400349",
    "2021-10-05 14:20:18.400351",
    "2021-10-05 14:20:18.400354",
    "2021-10-05 14:20:18.400356",
    "2021-10-05 14:20:18.400359",
    "2021-10-05 14:20:18.400363",
    "2021-10-05 14:20:18.400367",
    "2021-10-05 14:20:18.400369",
    "2021-10-05 14:20:18.400372",
    "2021-10-05 14:20:18.400375",
    "2021-10-05 14:20:18.400377",
    "2021-10-05 14:20:18.400380",
    "2021-10-05 14:20:18.400382",
    "2021-10-05 14:20:18.400384",
    "2021-10-05 14:20:18.400387",
    "2021-10-05 14:20:18.400393",
    "2021-10-05 14:20:18.400395",
    "2021-10-05 14:20:18.400398",
    "2021-10-05 14:20:18.400400",
    "2021-10-05 14:20:18.400403",
    "2021-10-05 14:20:18.400405",
    "2021-10-05 14:20:18.400408",
    "2021-10-05 14:20:18.400410",
    "2021-10-05 14:20:18.400413",
    "2021-10-05 14:20:18.400415",
    "2021-10-05 14:20:18.400419",
    "2021-10-05 14:20:18.400422"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# --- Add the "metadata" sheet after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$metaSheet.Name = "metadata"

# Header row (values first, then copy formatting from the styled "data" header cells)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Familial Meniere Disease"
$metaSheet.Range("C2").Value = 394
# data_version ("1.1") is stored as text in the source file, not a number --
# force text formatting first so COM doesn't silently coerce it to a double.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.1"
$metaSheet.Range("E2").Value = "2018-01-17T16:26:29.432517Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:18.396811"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/394/?format=json"

# Copy cell formatting (style index 1: bold, centered, bordered) from the
# "data" sheet header/index cells onto the corresponding "metadata" cells,
# without disturbing the values already written above.
$ws1.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)

$ws1.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Output "edit complete"
